# "Generate Report for Handback" - localization-status report refresh.
#
# The status for the de-de / zh-cn handback rows flips from "Ready for
# handoff" to "Handed back: in sync with en-US", the "Latest Handback
# DateTime" timestamps advance to the new handback run, and the stale
# "Error Detail" text (the old "handback file is not latest" warning) is
# cleared now that the handback is in sync.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- zh-cn sheet ---------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("K2").Value = "2016-08-13 05:00:18"
$wsZh.Range("P2").Value = ""

# --- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("K2").Value = "2016-08-13 05:00:28"
$wsDe.Range("P2").Value = ""

# --- Overview sheet: the per-language status columns mirror the same text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns("E:F").AutoFit()

$wsZh.Columns("C:C").AutoFit()
$wsZh.Columns("P:P").AutoFit()

$wsDe.Columns("C:C").AutoFit()
$wsDe.Columns("P:P").AutoFit()
